# Insert a new data row at row 225 (pushing the existing rows 225..257 down
# to 226..258) and populate it with a new price-record for "Ají" / "Inferno".
#
# The sheet is a flat price table: columns A..R, header in row 1, data rows
# starting at row 2. The new row reuses the constant columns (A, B, C, R)
# from the row that is being pushed down and fills in the new record's own
# values for D..Q.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 225:257 down to 226:258 by inserting a new blank row at 225.
$ws.Rows.Item(225).Insert()

# New row 225 values
$ws.Cells.Item(225, 1).Value = 8
$ws.Cells.Item(225, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(225, 3).Value = "Coquimbo"
$ws.Cells.Item(225, 4).Value = 44776
$ws.Cells.Item(225, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(225, 5).Value = 4
$ws.Cells.Item(225, 6).Value = 100112021
$ws.Cells.Item(225, 7).Value = "Ají"
$ws.Cells.Item(225, 8).Value = "Inferno"
$ws.Cells.Item(225, 9).Value = "Primera"
$ws.Cells.Item(225, 10).Value = 400
$ws.Cells.Item(225, 11).Value = 14000
$ws.Cells.Item(225, 12).Value = 15000
$ws.Cells.Item(225, 13).Value = 14500
$ws.Cells.Item(225, 14).Value = "`$/caja 12 kilos"
$ws.Cells.Item(225, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(225, 16).Value = 1208
$ws.Cells.Item(225, 17).Value = 12
$ws.Cells.Item(225, 18).Value = "Hortaliza"
